# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '51.635.71'
$ws.Range("E2").Value = '  -1.05%  '
# Row 3
$ws.Range("D3").Value = '2.920.34'
$ws.Range("E3").Value = '  +1.19%  '
# Row 4
$ws.Range("E4").Value = '  +0.05%  '
# Row 5
$ws.Range("D5").Value = '''354.08'
$ws.Range("E5").Value = '  +0.66%  '
# Row 6
$ws.Range("D6").Value = '''109.99'
$ws.Range("E6").Value = '  -1.57%  '
# Row 7
$ws.Range("D7").Value = '''0.560'
$ws.Range("E7").Value = '  -0.28%  '
# Row 8
$ws.Range("E8").Value = '  +0.04%  '
# Row 9
$ws.Range("D9").Value = '''0.636'
$ws.Range("E9").Value = '  +2.03%  '
# Row 10
$ws.Range("D10").Value = '''39.07'
$ws.Range("E10").Value = '  -2.82%  '
# Row 11
$ws.Range("D11").Value = '''0.0887'
$ws.Range("E11").Value = '  +3.00%  '
# Row 12
$ws.Range("E12").Value = '  +0.92%  '
# Row 13
$ws.Range("D13").Value = '''19.68'
$ws.Range("E13").Value = '  -2.08%  '
# Row 14
$ws.Range("D14").Value = '''7.95'
$ws.Range("E14").Value = '  +1.35%  '
# Row 15
$ws.Range("D15").Value = '3.383.98'
$ws.Range("E15").Value = '  +1.41%  '
# Row 16
$ws.Range("D16").Value = '2.938.56'
$ws.Range("E16").Value = '  +1.66%  '
# Row 17
$ws.Range("E17").Value = '  -1.44%  '
# Row 18
$ws.Range("D18").Value = '51.711.85'
$ws.Range("E18").Value = '  -0.85%  '
# Row 19
$ws.Range("D19").Value = '''7.59'
$ws.Range("E19").Value = '  -1.11%  '
# Row 20
$ws.Range("D20").Value = '''3.27'
$ws.Range("E20").Value = '  -2.78%  '
# Row 21
$ws.Range("D21").Value = '''14.19'
$ws.Range("E21").Value = '  +3.94%  '
# Row 22
$ws.Range("D22").Value = '0.0₃0980'
$ws.Range("E22").Value = '  -0.28%  '
# Row 23
$ws.Range("D23").Value = '''70.79'
$ws.Range("E23").Value = '  -0.45%  '
# Row 24
$ws.Range("D24").Value = '''269.58'
$ws.Range("E24").Value = '  +0.01%  '
# Row 25
$ws.Range("E25").Value = '  +1.19%  '
# Row 26
$ws.Range("D26").Value = '''0.183'
$ws.Range("E26").Value = '  +11.05%  '
# Row 27
$ws.Range("D27").Value = '''27.17'
$ws.Range("E27").Value = '  +2.76%  '
# Row 28
$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").Value = '''7.54'
$ws.Range("E28").Value = '  +18.91%  '
# Row 29
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  +0.10%  '
# Row 30
$ws.Range("E30").Value = '  +15.24%  '
# Row 31
$ws.Range("D31").Value = '''10.60'
$ws.Range("E31").Value = '  +0.07%  '
# Row 32
$ws.Range("D32").Value = '''37.71'
$ws.Range("E32").Value = '  -2.85%  '
# Row 33
$ws.Range("D33").Value = '''6.16'
$ws.Range("E33").Value = '  +3.15%  '
# Row 34
$ws.Range("D34").Value = '''52.28'
$ws.Range("E34").Value = '  -1.57%  '
# Row 35
$ws.Range("D35").Value = '''0.0442'
$ws.Range("E35").Value = '  -3.68%  '
# Row 36
$ws.Range("E36").Value = '  +0.06%  '
# Row 37
$ws.Range("E37").Value = '  -16.46%  '
# Row 38
$ws.Range("D38").Value = '''3.25'
$ws.Range("E38").Value = '  -1.87%  '
# Row 39
$ws.Range("E39").Value = '  -1.15%  '
# Row 40
$ws.Range("D40").Value = '''18.18'
$ws.Range("E40").Value = '  -2.62%  '
# Row 41
$ws.Range("E41").Value = '  +4.27%  '
# Row 42
$ws.Range("E42").Value = '  +0.88%  '
# Row 43
$ws.Range("D43").Value = '''23.05'
$ws.Range("E43").Value = '  +2.02%  '
# Row 44
$ws.Range("D44").Value = '''2.17'
$ws.Range("E44").Value = '  -2.17%  '
# Row 45
$ws.Range("E45").Value = '  +1.92%  '
# Row 46
$ws.Range("E46").Value = '  -4.76%  '
# Row 47
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.140.12'
$ws.Range("E47").Value = '  -2.23%  '
# Row 48
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '''3.44'
$ws.Range("E48").Value = '  -3.71%  '
# Row 49
$ws.Range("D49").Value = '''0.251'
$ws.Range("E49").Value = '  -2.52%  '
# Row 50
$ws.Range("D50").Value = '''0.0331'
$ws.Range("E50").Value = '  +2.53%  '
# Row 51
$ws.Range("D51").Value = '''9.15'
$ws.Range("E51").Value = '  +2.08%  '
